# Update the cryptocurrency price/volume table with the latest scraped values.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Cells in column D hold price strings (e.g. "64.207.05") that look numeric to Excel.
# Force them to Text format first so assigning the value keeps them as strings
# (matching the original inline-string cell type) instead of being parsed into numbers.
$dCells = @("D2","D3","D5","D6","D8","D9","D11","D14","D15","D16","D17","D20","D21","D22","D23","D24","D25","D28","D30","D34","D35","D36","D37","D38","D39","D41","D42","D45","D46","D49","D50")
foreach ($addr in $dCells) {
    $ws.Range($addr).NumberFormat = "@"
}

$ws.Range("D2").Value = '64.207.05'
$ws.Range("E2").Value = '  +3.55%  '

$ws.Range("D3").Value = '3.065.37'
$ws.Range("E3").Value = '  +2.38%  '

$ws.Range("E4").Value = '  +0.04%  '

$ws.Range("D5").Value = '561.23'
$ws.Range("E5").Value = '  +3.66%  '

$ws.Range("D6").Value = '143.56'
$ws.Range("E6").Value = '  +4.30%  '

$ws.Range("E7").Value = '  +0.06%  '

$ws.Range("D8").Value = '3.065.11'
$ws.Range("E8").Value = '  +2.35%  '

$ws.Range("D9").Value = '0.513'
$ws.Range("E9").Value = '  +5.47%  '

$ws.Range("E10").Value = '  +6.47%  '

$ws.Range("D11").Value = '6.11'
$ws.Range("E11").Value = '  -8.63%  '

$ws.Range("E12").Value = '  +10.43%  '

$ws.Range("E13").Value = '  +6.02%  '

$ws.Range("D14").Value = '35.60'
$ws.Range("E14").Value = '  +5.62%  '

$ws.Range("D15").Value = '3.564.11'
$ws.Range("E15").Value = '  +2.43%  '

$ws.Range("D16").Value = '64.210.60'
$ws.Range("E16").Value = '  +3.63%  '

$ws.Range("D17").Value = '3.064.79'
$ws.Range("E17").Value = '  +2.37%  '

$ws.Range("E18").Value = '  +2.87%  '

$ws.Range("E19").Value = '  +3.98%  '

$ws.Range("D20").Value = '479.72'
$ws.Range("E20").Value = '  +3.47%  '

$ws.Range("D21").Value = '14.00'
$ws.Range("E21").Value = '  +5.49%  '

$ws.Range("D22").Value = '0.683'
$ws.Range("E22").Value = '  +5.35%  '

$ws.Range("D23").Value = '7.57'
$ws.Range("E23").Value = '  +5.66%  '

$ws.Range("D24").Value = '14.33'
$ws.Range("E24").Value = '  +14.84%  '

$ws.Range("D25").Value = '82.32'
$ws.Range("E25").Value = '  +4.32%  '

$ws.Range("E26").Value = '  +0.07%  '

$ws.Range("E27").Value = '  +4.26%  '

$ws.Range("D28").Value = '8.10'
$ws.Range("E28").Value = '  +7.18%  '

$ws.Range("E29").Value = '  +2.59%  '

$ws.Range("D30").Value = '1.00'
$ws.Range("E30").Value = '  -0.06%  '

$ws.Range("E31").Value = '  +4.31%  '

$ws.Range("E32").Value = '  +2.27%  '

$ws.Range("E33").Value = '  +5.81%  '

$ws.Range("D34").Value = '5.74'
$ws.Range("E34").Value = '  +4.23%  '

$ws.Range("D35").Value = '6.28'
$ws.Range("E35").Value = '  +8.56%  '

$ws.Range("D36").Value = '55.01'
$ws.Range("E36").Value = '  +0.82%  '

$ws.Range("D37").Value = '0.0409'
$ws.Range("E37").Value = '  +6.22%  '

$ws.Range("D38").Value = '445.77'
$ws.Range("E38").Value = '  -0.57%  '

$ws.Range("D39").Value = '0.0815'
$ws.Range("E39").Value = '  +1.62%  '

$ws.Range("E40").Value = '  +12.83%  '

$ws.Range("D41").Value = '3.002.91'
$ws.Range("E41").Value = '  +2.79%  '

$ws.Range("D42").Value = '8.28'
$ws.Range("E42").Value = '  +3.30%  '

$ws.Range("E43").Value = '  +1.01%  '

$ws.Range("E44").Value = '  +4.85%  '

$ws.Range("D45").Value = '0.264'
$ws.Range("E45").Value = '  +7.59%  '

$ws.Range("D46").Value = '2.19'
$ws.Range("E46").Value = '  +10.38%  '

$ws.Range("E47").Value = '  -0.01%  '

$ws.Range("E48").Value = '  +4.66%  '

$ws.Range("D49").Value = '0.0₃0521'
$ws.Range("E49").Value = '  +5.89%  '

$ws.Range("D50").Value = '118.70'
$ws.Range("E50").Value = '  +3.56%  '

$ws.Range("E51").Value = '  +5.32%  '

# Restore the default (Normal) style on the price cells so no stray number-format
# style survives on the cells themselves (the text is already stored by this point).
foreach ($addr in $dCells) {
    $ws.Range($addr).Style = "Normal"
}
